# Generate Report for Handback
# Adds a new handback record (file 7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.md,
# "in sync with en-US") as a 4th data row to all three report sheets:
# Overview, zh-cn, de-de. Each sheet's table grows by one row, plus the
# corresponding hyperlinks are added.

$wb = $excel.ActiveWorkbook

$newFileName = "7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.md"
$newPathName = "e2e\7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.md"
$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $newFileName
$wsOverview.Range("B4").Value = $newPathName
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = "2016-08-24 22:45:23"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/997820e3283efca2f258d439bdaa9b0f5ed0a5b3/e2e/7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.md",
    [Type]::Missing,
    [Type]::Missing,
    $newPathName
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = $newFileName
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $statusInSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
# "True"/"False" look like booleans to the COM layer's auto-detection, so
# force literal text with a leading apostrophe, then strip the resulting
# quote-prefix style back to Normal to match the source workbook (which
# stores these as plain shared-string text, not boolean cells).
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("F4").Style = "Normal"
$wsZhCn.Range("G4").Value = "7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.997820e3283efca2f258d439bdaa9b0f5ed0a5b3.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-24 22:45:18"
$wsZhCn.Range("I4").Value = $newFileName
$wsZhCn.Range("J4").Value = "7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.997820e3283efca2f258d439bdaa9b0f5ed0a5b3.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-24 22:45:36"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("M4").Style = "Normal"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("O4").Style = "Normal"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/997820e3283efca2f258d439bdaa9b0f5ed0a5b3/e2e/7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.md",
    [Type]::Missing,
    [Type]::Missing,
    $newFileName
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/997820e3283efca2f258d439bdaa9b0f5ed0a5b3/e2e/7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.md",
    [Type]::Missing,
    [Type]::Missing,
    $newFileName
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = $newFileName
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $statusInSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
# see zh-cn comment above re: forcing literal "True"/"False" text
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("F4").Style = "Normal"
$wsDeDe.Range("G4").Value = "7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.997820e3283efca2f258d439bdaa9b0f5ed0a5b3.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-24 22:45:23"
$wsDeDe.Range("I4").Value = $newFileName
$wsDeDe.Range("J4").Value = "7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.997820e3283efca2f258d439bdaa9b0f5ed0a5b3.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-24 22:45:43"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("M4").Style = "Normal"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("O4").Style = "Normal"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/997820e3283efca2f258d439bdaa9b0f5ed0a5b3/e2e/7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.md",
    [Type]::Missing,
    [Type]::Missing,
    $newFileName
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/997820e3283efca2f258d439bdaa9b0f5ed0a5b3/e2e/7e08dc8d-bcd6-49ff-a5dc-e332dbced2ca.md",
    [Type]::Missing,
    [Type]::Missing,
    $newFileName
) | Out-Null

$wsOverview.Select()
$wsOverview.Range("A1").Select()
